$wb = $excel.ActiveWorkbook

# --- Sheet "Products": register the sale against the stocked product ---
$products = $wb.Worksheets.Item("Products")
$products.Range("F2").Value = 18
$products.Range("L2").Value = "2025-09-23T13:05:26.195Z"

# --- Sheet "Sales": append the new sale record as row 4 ---
$sales = $wb.Worksheets.Item("Sales")

$sales.Range("A4").Value = "07291cea-a90f-4352-b915-1121e63dcb59"
$sales.Range("B4").Value = 1

# Numeric-looking text must be forced to text; revert the style afterwards
# so the cell keeps the plain (unstyled) "number stored as text" look.
$sales.Range("C4").NumberFormat = "@"
$sales.Range("C4").Value = "500"
$sales.Range("C4").Style = "Normal"

# Empty-string text cells (as opposed to truly blank cells) round-trip via
# a bare quote-prefixed value, then drop the resulting quote-prefix style.
$sales.Range("D4").Value = "'"
$sales.Range("D4").Style = "Normal"

$sales.Range("E4").Value = "Cash"
$sales.Range("F4").Value = "Admin"

$sales.Range("G4").Value = "'"
$sales.Range("G4").Style = "Normal"

$sales.Range("H4").Value = "15ca83ea-a74c-421a-911c-b93c602bde13"
$sales.Range("I4").Value = "Amul Butter (500g)"

$sales.Range("J4").NumberFormat = "@"
$sales.Range("J4").Value = "500"
$sales.Range("J4").Style = "Normal"

$sales.Range("K4").NumberFormat = "@"
$sales.Range("K4").Value = "400"
$sales.Range("K4").Style = "Normal"

$sales.Range("L4").Value = "2025-09-23T13:05:26.171Z"
